$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "903.52.718"
$ws.Range("B2").Value = 42
$ws.Range("C2").Value = "stokta mevcut"
$ws.Range("D2").Value = "315,34 TL"
$ws.Range("E2").Value = "210,22 TL"
$ws.Range("F2").Value = "273,29 TL"

$ws.Range("A3").Value = "903.53.718"
$ws.Range("B3").Value = 21
$ws.Range("C3").Value = "stokta mevcut"
$ws.Range("D3").Value = "657,64 TL"
$ws.Range("E3").Value = "438,43 TL"
$ws.Range("F3").Value = "569,96 TL"

$ws.Range("A4").Value = "903.58.055"
$ws.Range("B4").Value = 87
$ws.Range("C4").Value = "stokta mevcut"
$ws.Range("D4").Value = "136,49 TL"
$ws.Range("E4").Value = "90,99 TL"
$ws.Range("F4").Value = "118,29 TL"

$ws.Range("A5").Value = "903.58.056"
$ws.Range("B5").Value = 55
$ws.Range("C5").Value = "stokta mevcut"
$ws.Range("D5").Value = "136,49 TL"
$ws.Range("E5").Value = "90,99 TL"
$ws.Range("F5").Value = "118,29 TL"

$ws.Range("A6").Value = "903.58.064"
$ws.Range("B6").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("C6").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("D6").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("E6").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("F6").Value = "urun hafele.com.tr de bulunmuyor"

$ws.Range("A7").Value = "903.58.057"
$ws.Range("B7").Value = 66
$ws.Range("C7").Value = "stokta mevcut"
$ws.Range("D7").Value = "383,61 TL"
$ws.Range("E7").Value = "255,74 TL"
$ws.Range("F7").Value = "332,46 TL"

$ws.Range("A8").Value = "903.58.068"
$ws.Range("B8").Value = 93
$ws.Range("C8").Value = "stokta mevcut"
$ws.Range("D8").Value = "148,03 TL"
$ws.Range("E8").Value = "98,69 TL"
$ws.Range("F8").Value = "128,29 TL"

$ws.Range("A9").Value = "903.58.070"
$ws.Range("B9").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("C9").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("D9").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("E9").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("F9").Value = "urun hafele.com.tr de bulunmuyor"

$ws.Range("A10").Value = "903.58.067"
$ws.Range("B10").Value = 847
$ws.Range("C10").Value = "stokta mevcut"
$ws.Range("D10").Value = "114,38 TL"
$ws.Range("E10").Value = "76,25 TL"
$ws.Range("F10").Value = "99,13 TL"

$ws.Range("A11").Value = "903.58.114"
$ws.Range("B11").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("C11").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("D11").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("E11").Value = "urun hafele.com.tr de bulunmuyor"
$ws.Range("F11").Value = "urun hafele.com.tr de bulunmuyor"

$ws.Range("A12").Value = "903.58.267"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "stokta mevcut"
$ws.Range("D12").Value = "72,06 TL"
$ws.Range("E12").Value = "48,04 TL"
$ws.Range("F12").Value = "62,46 TL"

$ws.Range("A13").Value = "903.58.323"
$ws.Range("B13").Value = 806
$ws.Range("C13").Value = "stokta mevcut"
$ws.Range("D13").Value = "157,64 TL"
$ws.Range("E13").Value = "105,10 TL"
$ws.Range("F13").Value = "136,63 TL"

$ws.Range("A14").Value = "903.58.368"
$ws.Range("B14").Value = 143
$ws.Range("C14").Value = "stokta mevcut"
$ws.Range("D14").Value = "99,95 TL"
$ws.Range("E14").Value = "66,63 TL"
$ws.Range("F14").Value = "86,63 TL"

$ws.Range("A15").Value = "903.70.124"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "3 ila 5 gün içinde"
$ws.Range("D15").Value = "614,24 TL"
$ws.Range("E15").Value = "409,49 TL"
$ws.Range("F15").Value = "573,29 TL"

Write-Output "done"